# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F8").Value = 2175
$ws1.Range("F12").Value = 1701
$ws1.Range("F17").Value = 251
$ws1.Range("F19").Value = 269
$ws1.Range("F20").Value = 1330
$ws1.Range("F22").Value = 278
$ws1.Range("F23").Value = 633
$ws1.Range("F24").Value = 12401
$ws1.Range("F25").Value = 12450
$ws1.Range("F26").Value = 923
$ws1.Range("F31").Value = 413
$ws1.Range("F32").Value = 1938

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 110

# --- Sheet "全部类型" (all types, aggregate) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F9").Value = 2175
$ws4.Range("F12").Value = 110
$ws4.Range("F14").Value = 1701
$ws4.Range("F22").Value = 251
$ws4.Range("F24").Value = 269
$ws4.Range("F25").Value = 1330
$ws4.Range("F27").Value = 278
$ws4.Range("F28").Value = 1
$ws4.Range("F29").Value = 633
$ws4.Range("F30").Value = 12401
$ws4.Range("F31").Value = 12450
$ws4.Range("F32").Value = 923
$ws4.Range("F37").Value = 413
$ws4.Range("F40").Value = 1938
